$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.296436
$ws.Range("H2").Value = 3.889308000000001
$ws.Range("I2").Value = 0.1944674516316147
$ws.Range("J2").Value = 0.1944674516316147
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 4.453045666666667
$ws.Range("N2").Value = 13.359137
$ws.Range("O2").Value = 0.4394129038053478
$ws.Range("P2").Value = 0.4394129038053478
$ws.Range("Q2").Value = 5.773088711910668
$ws.Range("R2").Value = 51.95779840719601
$ws.Range("S2").Value = 0.08545150761707386
$ws.Range("T2").Value = 0.08545150761707385

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.296436
$ws.Range("H3").Value = 3.889308000000001
$ws.Range("I3").Value = 0.1944674516316147
$ws.Range("J3").Value = 0.1944674516316147
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.150099666666667
$ws.Range("N3").Value = 9.450299000000001
$ws.Range("O3").Value = 0.3108421842981904
$ws.Range("P3").Value = 0.3108421842981904
$ws.Range("Q3").Value = 4.083902611454667
$ws.Range("R3").Value = 36.75512350309201
$ws.Range("S3").Value = 0.06044868744007383
$ws.Range("T3").Value = 0.06044868744007381

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.296436
$ws.Range("H4").Value = 3.889308000000001
$ws.Range("I4").Value = 0.1944674516316147
$ws.Range("J4").Value = 0.1944674516316147
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.530935
$ws.Range("N4").Value = 7.592805
$ws.Range("O4").Value = 0.2497449118964618
$ws.Range("P4").Value = 0.2497449118964618
$ws.Range("Q4").Value = 3.28119524766
$ws.Range("R4").Value = 29.53075722894
$ws.Range("S4").Value = 0.04856725657446708
$ws.Range("T4").Value = 0.04856725657446708

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.837275666666667
$ws.Range("H5").Value = 5.511827
$ws.Range("I5").Value = 0.2755942575194169
$ws.Range("J5").Value = 0.2755942575194169
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.453045666666667
$ws.Range("N5").Value = 13.359137
$ws.Range("O5").Value = 0.4394129038053478
$ws.Range("P5").Value = 0.4394129038053478
$ws.Range("Q5").Value = 8.181472445922113
$ws.Range("R5").Value = 73.633252013299
$ws.Range("S5").Value = 0.1210996729686858
$ws.Range("T5").Value = 0.1210996729686857

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.837275666666667
$ws.Range("H6").Value = 5.511827
$ws.Range("I6").Value = 0.2755942575194169
$ws.Range("J6").Value = 0.2755942575194169
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.150099666666667
$ws.Range("N6").Value = 9.450299000000001
$ws.Range("O6").Value = 0.3108421842981904
$ws.Range("P6").Value = 0.3108421842981904
$ws.Range("Q6").Value = 5.787601465141446
$ws.Range("R6").Value = 52.08841318627301
$ws.Range("S6").Value = 0.08566632098737353
$ws.Range("T6").Value = 0.08566632098737352

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.837275666666667
$ws.Range("H7").Value = 5.511827
$ws.Range("I7").Value = 0.2755942575194169
$ws.Range("J7").Value = 0.2755942575194169
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.530935
$ws.Range("N7").Value = 7.592805
$ws.Range("O7").Value = 0.2497449118964618
$ws.Range("P7").Value = 0.2497449118964618
$ws.Range("Q7").Value = 4.650025289415
$ws.Range("R7").Value = 41.850227604735
$ws.Range("S7").Value = 0.06882826356335757
$ws.Range("T7").Value = 0.06882826356335757

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.532884666666666
$ws.Range("H8").Value = 10.598654
$ws.Range("I8").Value = 0.5299382908489685
$ws.Range("J8").Value = 0.5299382908489685
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 4.453045666666667
$ws.Range("N8").Value = 13.359137
$ws.Range("O8").Value = 0.4394129038053478
$ws.Range("P8").Value = 0.4394129038053478
$ws.Range("Q8").Value = 15.73209675573311
$ws.Range("R8").Value = 141.588870801598
$ws.Range("S8").Value = 0.2328617232195882
$ws.Range("T8").Value = 0.2328617232195882

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.532884666666666
$ws.Range("H9").Value = 10.598654
$ws.Range("I9").Value = 0.5299382908489685
$ws.Range("J9").Value = 0.5299382908489685
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.150099666666667
$ws.Range("N9").Value = 9.450299000000001
$ws.Range("O9").Value = 0.3108421842981904
$ws.Range("P9").Value = 0.3108421842981904
$ws.Range("Q9").Value = 11.12893881083845
$ws.Range("R9").Value = 100.160449297546
$ws.Range("S9").Value = 0.1647271758707431
$ws.Range("T9").Value = 0.1647271758707431

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.532884666666666
$ws.Range("H10").Value = 10.598654
$ws.Range("I10").Value = 0.5299382908489685
$ws.Range("J10").Value = 0.5299382908489685
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.530935
$ws.Range("N10").Value = 7.592805
$ws.Range("O10").Value = 0.2497449118964618
$ws.Range("P10").Value = 0.2497449118964618
$ws.Range("Q10").Value = 8.94150145383
$ws.Range("R10").Value = 80.47351308447
$ws.Range("S10").Value = 0.1323493917586372
$ws.Range("T10").Value = 0.1323493917586372
